$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.525.31"
$ws.Range("E2").Value = "  -2.65%  "

$ws.Range("D3").Value = "'1.812.72"
$ws.Range("E3").Value = "  -2.17%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.48%  "

$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").Value = "'308.12"
$ws.Range("E6").Value = "  -1.91%  "

$ws.Range("D7").Value = "'0.4558"
$ws.Range("E7").Value = "  -2.00%  "

$ws.Range("D8").Value = "'0.3666"
$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("D9").Value = "'0.07136"
$ws.Range("E9").Value = "  -2.16%  "

$ws.Range("D10").Value = "'0.8786"

$ws.Range("D11").Value = "'0.07771"
$ws.Range("E11").Value = "  -1.14%  "

$ws.Range("D12").Value = "'19.39"
$ws.Range("E12").Value = "  -3.64%  "

$ws.Range("D13").Value = "'1.797.23"
$ws.Range("E13").Value = "  -2.73%  "

$ws.Range("D14").Value = "'5.289"
$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("D15").Value = "'6.362"
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("D16").Value = "'86.42"
$ws.Range("E16").Value = "  -5.25%  "

$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").Value = "'0.000008606"
$ws.Range("E18").Value = "  -3.53%  "

$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").Value = "'26.591.39"
$ws.Range("E20").Value = "  -2.50%  "

$ws.Range("E21").Value = "  -3.12%  "

$ws.Range("D22").Value = "'5.006"
$ws.Range("E22").Value = "  -1.48%  "

$ws.Range("D23").Value = "'10.47"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").Value = "'1.984"
$ws.Range("E24").Value = "  +1.58%  "

$ws.Range("D25").Value = "'151.69"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("D26").Value = "'17.95"
$ws.Range("E26").Value = "  -2.42%  "

$ws.Range("D27").Value = "'2.058"
$ws.Range("E27").Value = "  +0.80%  "

$ws.Range("D28").Value = "'112.89"
$ws.Range("E28").Value = "  -2.63%  "

$ws.Range("D29").Value = "'4.859"
$ws.Range("E29").Value = "  -3.71%  "

$ws.Range("D30").Value = "'0.08685"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("D31").Value = "'3.064"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("D32").Value = "'4.508"
$ws.Range("E32").Value = "  -0.24%  "

$ws.Range("D33").Value = "'0.7342"
$ws.Range("E33").Value = "  -4.14%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'2.695"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.120"
$ws.Range("E35").Value = "  -4.09%  "

$ws.Range("E36").Value = "  +0.50%  "

$ws.Range("D37").Value = "'1.084"
$ws.Range("E37").Value = "  -2.59%  "

$ws.Range("D38").Value = "'0.01952"
$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("D39").Value = "'0.05119"
$ws.Range("E39").Value = "  -1.99%  "

$ws.Range("D40").Value = "'2.906"
$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("D41").Value = "'6.995"
$ws.Range("E41").Value = "  -0.72%  "

$ws.Range("D42").Value = "'0.5014"
$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("E43").Value = "  -4.13%  "

$ws.Range("D44").Value = "'8.177"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("D46").Value = "'0.4613"
$ws.Range("E46").Value = "  -3.85%  "

$ws.Range("D47").Value = "'9.959"
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("D48").Value = "'100.94"
$ws.Range("E48").Value = "  -1.66%  "

$ws.Range("D49").Value = "'1.592"
$ws.Range("E49").Value = "  -3.13%  "

$ws.Range("D50").Value = "'0.05999"
$ws.Range("E50").Value = "  -3.28%  "

$ws.Range("E51").Value = "  -1.44%  "
